$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price value is a "plain" number (single decimal point) need to be
# forced to Text format first, so Excel stores them as text (matching the original
# inline-string cells) instead of auto-converting them to numeric values.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '51.710.73'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').Value = '2.947.42'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '380.13'
$ws.Range('E5').Value = '  +7.52%  '
$ws.Range('D6').Value = '105.29'
$ws.Range('E6').Value = '  -0.62%  '
$ws.Range('E7').Value = '  -1.42%  '
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('E10').Value = '  -0.67%  '
$ws.Range('E11').Value = '  -0.73%  '
$ws.Range('E12').Value = '  -0.33%  '
$ws.Range('D13').Value = '18.59'
$ws.Range('E13').Value = '  -1.46%  '
$ws.Range('D14').Value = '3.407.93'
$ws.Range('E14').Value = '  -0.79%  '
$ws.Range('E15').Value = '  -0.14%  '
$ws.Range('D16').Value = '2.947.29'
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('D17').Value = '0.956'
$ws.Range('E17').Value = '  -2.79%  '
$ws.Range('D18').Value = '51.648.22'
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('D19').Value = '3.43'
$ws.Range('E19').Value = '  +3.42%  '
$ws.Range('D20').Value = '7.43'
$ws.Range('E20').Value = '  +1.34%  '
$ws.Range('D21').Value = '13.21'
$ws.Range('E21').Value = '  -0.70%  '
$ws.Range('D22').Value = '0.0₃0957'
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('D23').Value = '68.57'
$ws.Range('E23').Value = '  -0.48%  '
$ws.Range('D24').Value = '263.32'
$ws.Range('E24').Value = '  -0.83%  '
$ws.Range('E25').Value = '  +4.80%  '
$ws.Range('D26').Value = '7.47'
$ws.Range('E26').Value = '  +19.43%  '
$ws.Range('D27').Value = '4.16'
$ws.Range('E27').Value = '  -4.02%  '
$ws.Range('E28').Value = '  -3.70%  '
$ws.Range('D29').Value = '7.45'
$ws.Range('E29').Value = '  +2.80%  '
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('D31').Value = '25.97'
$ws.Range('E31').Value = '  -2.35%  '
$ws.Range('D32').Value = '0.104'
$ws.Range('E32').Value = '  -5.17%  '
$ws.Range('E33').Value = '  -1.32%  '
$ws.Range('D34').Value = '52.38'
$ws.Range('E34').Value = '  +3.12%  '
$ws.Range('D35').Value = '34.43'
$ws.Range('E35').Value = '  -2.88%  '
$ws.Range('E36').Value = '  -4.29%  '
$ws.Range('E37').Value = '  +2.12%  '
$ws.Range('E38').Value = '  +0.33%  '
$ws.Range('E39').Value = '  -5.44%  '
$ws.Range('D40').Value = '17.36'
$ws.Range('E40').Value = '  +0.44%  '
$ws.Range('D41').Value = '2.65'
$ws.Range('E41').Value = '  -5.74%  '
$ws.Range('E42').Value = '  -3.73%  '
$ws.Range('E43').Value = '  +0.10%  '
$ws.Range('D44').Value = '124.76'
$ws.Range('E44').Value = '  +1.04%  '
$ws.Range('D45').Value = '21.89'
$ws.Range('E45').Value = '  -4.06%  '
$ws.Range('D46').Value = '2.07'
$ws.Range('E46').Value = '  -4.02%  '
$ws.Range('D47').Value = '0.278'
$ws.Range('E47').Value = '  +17.35%  '
$ws.Range('D48').Value = '2.031.22'
$ws.Range('E48').Value = '  -3.53%  '
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('E50').Value = '  -1.06%  '
$ws.Range('D51').Value = '0.0327'
$ws.Range('E51').Value = '  +0.62%  '
